$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# A leading apostrophe forces Excel to keep numeric-looking strings (e.g. '98.65')
# stored as text, matching the original inline-string cell type.

$ws.Range("D2").Value = '43.771.40'
$ws.Range("E2").Value = '  +0.13%  '
$ws.Range("D3").Value = '2.297.38'
$ws.Range("E3").Value = '  -1.46%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = "'98.65"
$ws.Range("E5").Value = '  +3.32%  '
$ws.Range("D6").Value = "'270.78"
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").Value = "'0.626"
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").Value = "'0.607"
$ws.Range("E9").Value = '  -2.37%  '
$ws.Range("D10").Value = "'45.39"
$ws.Range("E10").Value = '  +0.21%  '
$ws.Range("D11").Value = "'0.0932"
$ws.Range("E11").Value = '  -1.22%  '
$ws.Range("D12").Value = "'7.89"
$ws.Range("E12").Value = '  -2.90%  '
$ws.Range("E13").Value = '  +1.68%  '
$ws.Range("D14").Value = "'15.88"
$ws.Range("E14").Value = '  +1.58%  '
$ws.Range("D15").Value = '2.641.60'
$ws.Range("E15").Value = '  -1.43%  '
$ws.Range("D16").Value = "'0.857"
$ws.Range("E16").Value = '  -0.54%  '
$ws.Range("D17").Value = '2.293.37'
$ws.Range("E17").Value = '  -1.73%  '
$ws.Range("D18").Value = '43.774.45'
$ws.Range("E18").Value = '  +0.26%  '
$ws.Range("E19").Value = '  +2.23%  '
$ws.Range("D20").Value = "'6.22"
$ws.Range("E20").Value = '  -3.06%  '
$ws.Range("E21").Value = '  -0.30%  '
$ws.Range("E22").Value = '  +7.78%  '
$ws.Range("D23").Value = "'233.38"
$ws.Range("E23").Value = '  -2.44%  '
$ws.Range("D24").Value = "'2.84"
$ws.Range("E24").Value = '  +12.35%  '
$ws.Range("D25").Value = "'9.14"
$ws.Range("E25").Value = '  -2.18%  '
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("D27").Value = "'11.31"
$ws.Range("E27").Value = '  -1.09%  '
$ws.Range("D28").Value = "'3.45"
$ws.Range("E28").Value = '  -1.07%  '
$ws.Range("D29").Value = "'2.29"
$ws.Range("E29").Value = '  +0.61%  '
$ws.Range("E30").Value = '  -0.23%  '
$ws.Range("D31").Value = "'176.67"
$ws.Range("E31").Value = '  +2.30%  '
$ws.Range("D32").Value = "'21.87"
$ws.Range("E32").Value = '  -3.47%  '
$ws.Range("D33").Value = "'0.0896"
$ws.Range("E33").Value = '  -0.31%  '
$ws.Range("E34").Value = '  -0.98%  '
$ws.Range("D35").Value = "'0.126"
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("D36").Value = "'4.75"
$ws.Range("E36").Value = '  +8.22%  '
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("D38").Value = "'0.0352"
$ws.Range("E38").Value = '  -2.91%  '
$ws.Range("E39").Value = '  +4.25%  '
$ws.Range("E40").Value = '  +0.31%  '
$ws.Range("D41").Value = "'2.34"
$ws.Range("E41").Value = '  -1.74%  '
$ws.Range("E42").Value = '  +0.97%  '
$ws.Range("D43").Value = "'12.20"
$ws.Range("E43").Value = '  +0.95%  '
$ws.Range("D44").Value = "'64.86"
$ws.Range("E44").Value = '  +4.26%  '
$ws.Range("D45").Value = "'8.85"
$ws.Range("E45").Value = '  -3.27%  '
$ws.Range("E46").Value = '  -2.52%  '
$ws.Range("E47").Value = '  -1.27%  '
$ws.Range("E48").Value = '  +0.62%  '
$ws.Range("D49").Value = "'98.05"
$ws.Range("E49").Value = '  -2.31%  '
$ws.Range("D50").Value = "'0.444"
$ws.Range("E50").Value = '  +6.82%  '
$ws.Range("D51").Value = "'1.53"
$ws.Range("E51").Value = '  +11.41%  '
